# Add a new "临床诊断病例" (clinically diagnosed cases) column of data (column E)
# to the Huanggang 2020-02-13 COVID-19 case-count sheet, and move the
# active selection from D13 to F13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column header (E2) ---------------------------------------------
$ws.Range("E2").Value = "临床诊断病例"

# --- New column data (E3:E12), one row per county/district --------------
$ws.Range("E3").Value  = -1   # 黄州
$ws.Range("E4").Value  = 0    # 团风
$ws.Range("E5").Value  = -9   # 红安
$ws.Range("E6").Value  = 0    # 罗田
$ws.Range("E7").Value  = 0    # 英山
$ws.Range("E8").Value  = -1   # 浠水
$ws.Range("E9").Value  = -8   # 蕲春
$ws.Range("E10").Value = 0    # 黄梅
$ws.Range("E11").Value = 0    # 麻城
$ws.Range("E12").Value = -15  # 武穴

# --- Totals row (确诊病例 / running total row) ---------------------------
$ws.Range("E13").Value = -34

# --- Sum row (全市累计) - formula mirroring the B/C/D sum formulas -------
$ws.Range("E14").Formula = "=SUM(E3:E12)"

# --- Formatting: match the integer number format used by the rest of the
#     table (columns B:D use the same "0" / numFmtId 1 style) -----------
$ws.Range("E2:E14").NumberFormat = "0"

# --- Move the active selection to F13 (as in the edited workbook) -------
$ws.Range("F13").Select()
